$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update: rows 96-98 get new prices for the current week, and the
# data that used to live in rows 96-98 (previous week) is preserved by
# appending it as new rows 99-101 at the bottom of the sheet.

$oldRows = @(96, 97, 98)
$newRowStart = 99

foreach ($r in $oldRows) {
    $destRow = $newRowStart + ($r - 96)

    # Copy column by column so that values are preserved exactly as they
    # were before this week's update.
    for ($col = 1; $col -le 18; $col++) {
        $srcCell = $ws.Cells.Item($r, $col)
        $dstCell = $ws.Cells.Item($destRow, $col)
        $dstCell.Value2 = $srcCell.Value2
    }

    # Column D (Fecha) carries a special date number format in the source
    # data; make sure the copied row keeps it too.
    $ws.Cells.Item($destRow, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat
}

# --- Update rows 96-98 with this week's new values.

# Row 96: Banquete
$ws.Range("D96").Value2 = 44516
$ws.Range("J96").Value2 = 250
$ws.Range("K96").Value2 = 1500
$ws.Range("L96").Value2 = 1500
$ws.Range("M96").Value2 = 1500
$ws.Range("O96").Value2 = "Provincia de Linares"
$ws.Range("P96").Value2 = 1500

# Row 97: Primera
$ws.Range("D97").Value2 = 44516
$ws.Range("J97").Value2 = 340
$ws.Range("K97").Value2 = 1300
$ws.Range("L97").Value2 = 1300
$ws.Range("M97").Value2 = 1300
$ws.Range("O97").Value2 = "Provincia de Linares"
$ws.Range("P97").Value2 = 1300

# Row 98: Segunda
$ws.Range("D98").Value2 = 44516
$ws.Range("J98").Value2 = 106
$ws.Range("K98").Value2 = 1100
$ws.Range("L98").Value2 = 1100
$ws.Range("M98").Value2 = 1100
$ws.Range("O98").Value2 = "Provincia de Linares"
$ws.Range("P98").Value2 = 1100
